# Update puerto_rico_stoch scenarios workbook for rivanna.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SolverSettings")

# Row 10: flip the "N" flags to "Y" for columns C through AA (B10 is already "Y").
$ws.Range("C10:AA10").Value = "Y"

# Update the view: scroll so column F is the left-most visible column and
# select cell AB10 (just past the used range).
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = $ws.Range("F1").Column
$ws.Range("AB10").Select()
